# "added colors to rows"
#
# Employee DTR (daily time record) sheet: highlight the daily rows to call
# out attendance status.
#   - Days with no Time-In/Time-Out AND no Overtime/Vacation-Leave entered
#     are flagged as a full day of Sick Leave: SICK LEAVE (col I) = 1 and
#     the row is shaded red.
#   - Days that already carry Overtime or Vacation-Leave hours (i.e. there
#     is something noteworthy about the day) are shaded blue so they stand
#     out from the plain rows.
# Also: row 19 col B is normalized from a stray space-string to a literal
# FALSE, and the FLOOR(...,1,1) formulas are simplified to the 2-arg
# FLOOR(...,1) form (the 3rd "mode" argument was a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel VBA/COM colors are BGR-packed integers (same as the RGB() macro).
$Blue   = 13411113   # FF29A3CC
$Orange = 6737151    # FFFFCC66 (kept in the palette, unused on this sheet)
$Red    = 6184671    # FFDF5E5E

$RedRows  = @(5, 6, 7, 8, 11, 12)
$BlueRows = @(9, 13, 14, 15)

foreach ($row in $RedRows) {
    $ws.Range("A" + $row + ":J" + $row).Interior.Color = $Red
    $ws.Range("I" + $row).Value = 1
}

foreach ($row in $BlueRows) {
    $ws.Range("A" + $row + ":J" + $row).Interior.Color = $Blue
}

# Row 19 is merged (A19:G19); the stored value under the merge's hidden
# cell B19 flips from a blank/space string to a literal boolean FALSE.
# A plain Range.Value assignment is a no-op on a non-anchor cell of a
# merged range, so use an array-formula write (which is permitted) to
# land the FALSE into B19.
$ws.Range("B19").FormulaArray = "=FALSE()"

# Simplify FLOOR(number, significance, mode) -> FLOOR(number, significance);
# the trailing "1" mode argument was redundant.
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
